# Applies cryptocurrency price/volume updates per commit:
# "Updated cryptos list on Mon Feb 20 20:46:25 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + "24.858.11"
$ws.Range("E2").Value = "  +1.15%  "
$ws.Range("D3").Value = "'" + "1.707.70"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("D4").Value = "'" + "1.002"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'" + "314.87"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").Value = "'" + "1.002"
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("D7").Value = "'" + "0.4019"
$ws.Range("E7").Value = "  +3.23%  "
$ws.Range("D8").Value = "'" + "0.4045"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'" + "1.471"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'" + "53.66"
$ws.Range("E11").Value = "  +1.68%  "
$ws.Range("D12").Value = "'" + "0.08796"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "'" + "26.33"
$ws.Range("E13").Value = "  +7.03%  "
$ws.Range("D14").Value = "'" + "7.507"
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").Value = "'" + "8.003"
$ws.Range("E15").Value = "  +0.86%  "
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "'" + "1.639.07"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "'" + "95.50"
$ws.Range("E18").Value = "  -2.43%  "
$ws.Range("D19").Value = "'" + "0.07180"
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "'" + "20.96"
$ws.Range("E20").Value = "  +6.77%  "
$ws.Range("D21").Value = "'" + "7.289"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'" + "1.003"
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "'" + "14.45"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").Value = "'" + "24.857.49"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "'" + "2.348"
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -3.57%  "
$ws.Range("D27").Value = "'" + "6.394"
$ws.Range("E27").Value = "  +22.37%  "
$ws.Range("D28").Value = "'" + "23.08"
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("D29").Value = "'" + "161.56"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").Value = "'" + "144.04"
$ws.Range("E30").Value = "  +6.02%  "
$ws.Range("D31").Value = "'" + "8.347"
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'" + "2.287"
$ws.Range("E32").Value = "  +15.69%  "
$ws.Range("D33").Value = "'" + "1.850.28"
$ws.Range("E33").Value = "  -0.91%  "
$ws.Range("D34").Value = "'" + "0.08712"
$ws.Range("E34").Value = "  -0.42%  "
$ws.Range("E35").Value = "  +10.02%  "
$ws.Range("D36").Value = "'" + "7.215"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("D37").Value = "'" + "1.027"
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("D38").Value = "'" + "0.2868"
$ws.Range("E38").Value = "  +5.98%  "
$ws.Range("D39").Value = "'" + "0.8416"
$ws.Range("E39").Value = "  +8.91%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'" + "0.09439"
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("D42").Value = "'" + "14.21"
$ws.Range("E42").Value = "  +1.69%  "
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("D44").Value = "'" + "17.58"
$ws.Range("E44").Value = "  +5.76%  "
$ws.Range("D45").Value = "'" + "2.718"
$ws.Range("E45").Value = "  +6.10%  "
$ws.Range("D46").Value = "'" + "0.7426"
$ws.Range("E46").Value = "  +4.37%  "
$ws.Range("D47").Value = "'" + "4.227"
$ws.Range("E47").Value = "  +0.68%  "
$ws.Range("D48").Value = "'" + "1.372"
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("D49").Value = "'" + "1.002"
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").Value = "'" + "140.51"
$ws.Range("E50").Value = "  +1.78%  "
$ws.Range("D51").Value = "'" + "0.08395"
$ws.Range("E51").Value = "  +5.54%  "
